$d = $word.ActiveDocument

# --- Paragraph 1: "purpose of the project" paragraph ---
$target1 = "purpose of the project"
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($target1)) {
        $p1 = $p
        break
    }
}
$r1 = $p1.Range
$r1b = $d.Range($r1.Start, $r1.End - 1)
$r1b.Text = "The purpose of the project is to create a WhatsApp/Slack/Discord style chat application. Users will sign into the application using a username and password. Users can search for other users by username, and add them to their “friends” list. Users can create forums based on different topics (called tags). Other users can comment on different forums. User will have admin privileges over the forum they create, and can create and assign different roles with varying privileges to other users. Admins of a particular forum can add or remove tags. Users can search for forums based on topics or tags. This project may in the future include additional features to expand its functionality."

# --- Paragraph 2: "application will include several different roles" paragraph ---
$target2 = "several different roles"
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($target2)) {
        $p2 = $p
        break
    }
}
$r2 = $p2.Range
$r2b = $d.Range($r2.Start, $r2.End - 1)
$r2b.Text = "The application will include several different roles. The user is the base role the application will support. The user will have the ability to search for other users, chat to other users, search for forums based on tags, and comment on forums. Once the user creates a forum, they gain the status of admin over their forum. The admin is an extension of the user, and is assigned to a user that creates a forum. The admin will have admin privileges ONLY over their specific forum. Admins can assign moderator or admin privileges to other users. Admins will assign tags to their forum on creation, and can edit these tags anytime. Super admins oversee the entire application, and have privileges that are distinct from all other users. Super admins can monitor any and all forums, as well as individual user chats. Super admins can assist users with technical issues, as well as modify or delete forums based on technical issues or infringement of user policies. "
